# Auto-generated edit script: updates TPM-derived NATMI metrics for Hras-Cav1 LR pairs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 4.455765
$ws.Cells.Item(2, 8).Value = 13.367295
$ws.Cells.Item(2, 9).Value = 0.1558824083674925
$ws.Cells.Item(2, 10).Value = 0.167793131187596
$ws.Cells.Item(2, 13).Value = 556.7425436666666
$ws.Cells.Item(2, 14).Value = 1670.227631
$ws.Cells.Item(2, 15).Value = 0.7235863858022448
$ws.Cells.Item(2, 16).Value = 0.7685368570853349
$ws.Cells.Item(2, 17).Value = 2480.713940080905
$ws.Cells.Item(2, 18).Value = 22326.42546072814
$ws.Cells.Item(2, 19).Value = 0.1127943884807835
$ws.Cells.Item(2, 20).Value = 0.1289552056834223

$ws.Cells.Item(3, 7).Value = 4.455765
$ws.Cells.Item(3, 8).Value = 13.367295
$ws.Cells.Item(3, 9).Value = 0.1558824083674925
$ws.Cells.Item(3, 10).Value = 0.167793131187596
$ws.Cells.Item(3, 15).Value = 0.0986194207087145
$ws.Cells.Item(3, 16).Value = 0.1047458342586422
$ws.Cells.Item(3, 17).Value = 338.102784291015
$ws.Cells.Item(3, 18).Value = 3042.925058619135
$ws.Cells.Item(3, 19).Value = 0.01537303281188138
$ws.Cells.Item(3, 20).Value = 0.01757563150911454

$ws.Cells.Item(4, 7).Value = 4.455765
$ws.Cells.Item(4, 8).Value = 13.367295
$ws.Cells.Item(4, 9).Value = 0.1558824083674925
$ws.Cells.Item(4, 10).Value = 0.167793131187596
$ws.Cells.Item(4, 13).Value = 1.402487333333333
$ws.Cells.Item(4, 14).Value = 4.207462
$ws.Cells.Item(4, 15).Value = 0.001822782814434402
$ws.Cells.Item(4, 16).Value = 0.001936017319896666
$ws.Cells.Item(4, 17).Value = 6.249153972809999
$ws.Cells.Item(4, 18).Value = 56.24238575528999
$ws.Cells.Item(4, 19).Value = 0.0002841397750449108
$ws.Cells.Item(4, 20).Value = 0.0003248504081388793

$ws.Cells.Item(5, 7).Value = 4.455765
$ws.Cells.Item(5, 8).Value = 13.367295
$ws.Cells.Item(5, 9).Value = 0.1558824083674925
$ws.Cells.Item(5, 10).Value = 0.167793131187596
$ws.Cells.Item(5, 13).Value = 135.006546
$ws.Cells.Item(5, 14).Value = 270.013092
$ws.Cells.Item(5, 15).Value = 0.1754651225976237
$ws.Cells.Item(5, 16).Value = 0.1242435517446983
$ws.Cells.Item(5, 17).Value = 601.55744243769
$ws.Cells.Item(5, 18).Value = 3609.34465462614
$ws.Cells.Item(5, 19).Value = 0.02735192589501491
$ws.Cells.Item(5, 20).Value = 0.02084721457711105

$ws.Cells.Item(6, 7).Value = 4.455765
$ws.Cells.Item(6, 8).Value = 13.367295
$ws.Cells.Item(6, 9).Value = 0.1558824083674925
$ws.Cells.Item(6, 10).Value = 0.167793131187596
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.3895486666666667
$ws.Cells.Item(6, 14).Value = 1.168646
$ws.Cells.Item(6, 15).Value = 0.0005062880769826339
$ws.Cells.Item(6, 16).Value = 0.0005377395914277917
$ws.Cells.Item(6, 17).Value = 1.73573731473
$ws.Cells.Item(6, 18).Value = 15.62163583257
$ws.Cells.Item(6, 19).Value = 0.00007892140476779942
$ws.Cells.Item(6, 20).Value = 0.00009022900980920775

$ws.Cells.Item(7, 9).Value = 0.1858758098371279
$ws.Cells.Item(7, 10).Value = 0.2000782799754709
$ws.Cells.Item(7, 13).Value = 556.7425436666666
$ws.Cells.Item(7, 14).Value = 1670.227631
$ws.Cells.Item(7, 15).Value = 0.7235863858022448
$ws.Cells.Item(7, 16).Value = 0.7685368570853349
$ws.Cells.Item(7, 17).Value = 2958.029179917062
$ws.Cells.Item(7, 18).Value = 26622.26261925356
$ws.Cells.Item(7, 19).Value = 0.1344972054481127
$ws.Cells.Item(7, 20).Value = 0.1537675324633881

$ws.Cells.Item(8, 9).Value = 0.1858758098371279
$ws.Cells.Item(8, 10).Value = 0.2000782799754709
$ws.Cells.Item(8, 15).Value = 0.0986194207087145
$ws.Cells.Item(8, 16).Value = 0.1047458342586422
$ws.Cells.Item(8, 19).Value = 0.01833096468990073
$ws.Cells.Item(8, 20).Value = 0.02095736635306488

$ws.Cells.Item(9, 9).Value = 0.1858758098371279
$ws.Cells.Item(9, 10).Value = 0.2000782799754709
$ws.Cells.Item(9, 13).Value = 1.402487333333333
$ws.Cells.Item(9, 14).Value = 4.207462
$ws.Cells.Item(9, 15).Value = 0.001822782814434402
$ws.Cells.Item(9, 16).Value = 0.001936017319896666
$ws.Cells.Item(9, 17).Value = 7.451556385724889
$ws.Cells.Item(9, 18).Value = 67.06400747152399
$ws.Cells.Item(9, 19).Value = 0.0003388112317901938
$ws.Cells.Item(9, 20).Value = 0.0003873550153676459

$ws.Cells.Item(10, 9).Value = 0.1858758098371279
$ws.Cells.Item(10, 10).Value = 0.2000782799754709
$ws.Cells.Item(10, 13).Value = 135.006546
$ws.Cells.Item(10, 14).Value = 270.013092
$ws.Cells.Item(10, 15).Value = 0.1754651225976237
$ws.Cells.Item(10, 16).Value = 0.1242435517446983
$ws.Cells.Item(10, 17).Value = 717.3033695569642
$ws.Cells.Item(10, 18).Value = 4303.820217341785
$ws.Cells.Item(10, 19).Value = 0.03261472176100423
$ws.Cells.Item(10, 20).Value = 0.02485843613112266

$ws.Cells.Item(11, 9).Value = 0.1858758098371279
$ws.Cells.Item(11, 10).Value = 0.2000782799754709
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.3895486666666667
$ws.Cells.Item(11, 14).Value = 1.168646
$ws.Cells.Item(11, 15).Value = 0.0005062880769826339
$ws.Cells.Item(11, 16).Value = 0.0005377395914277917
$ws.Cells.Item(11, 17).Value = 2.069711280565778
$ws.Cells.Item(11, 18).Value = 18.627401525092
$ws.Cells.Item(11, 19).Value = 0.00009410670632002924
$ws.Cells.Item(11, 20).Value = 0.000107590012527585

$ws.Cells.Item(12, 7).Value = 5.633732333333334
$ws.Cells.Item(12, 8).Value = 16.901197
$ws.Cells.Item(12, 9).Value = 0.1970929266282699
$ws.Cells.Item(12, 10).Value = 0.2121524785267629
$ws.Cells.Item(12, 13).Value = 556.7425436666666
$ws.Cells.Item(12, 14).Value = 1670.227631
$ws.Cells.Item(12, 15).Value = 0.7235863858022448
$ws.Cells.Item(12, 16).Value = 0.7685368570853349
$ws.Cells.Item(12, 17).Value = 3136.538469597145
$ws.Cells.Item(12, 18).Value = 28228.84622637431
$ws.Cells.Item(12, 19).Value = 0.1426137584461369
$ws.Cells.Item(12, 20).Value = 0.1630469990698224

$ws.Cells.Item(13, 7).Value = 5.633732333333334
$ws.Cells.Item(13, 8).Value = 16.901197
$ws.Cells.Item(13, 9).Value = 0.1970929266282699
$ws.Cells.Item(13, 10).Value = 0.2121524785267629
$ws.Cells.Item(13, 15).Value = 0.0986194207087145
$ws.Cells.Item(13, 16).Value = 0.1047458342586422
$ws.Cells.Item(13, 17).Value = 427.4867700272157
$ws.Cells.Item(13, 18).Value = 3847.380930244941
$ws.Cells.Item(13, 19).Value = 0.01943719024986515
$ws.Cells.Item(13, 20).Value = 0.02222208835332444

$ws.Cells.Item(14, 7).Value = 5.633732333333334
$ws.Cells.Item(14, 8).Value = 16.901197
$ws.Cells.Item(14, 9).Value = 0.1970929266282699
$ws.Cells.Item(14, 10).Value = 0.2121524785267629
$ws.Cells.Item(14, 13).Value = 1.402487333333333
$ws.Cells.Item(14, 14).Value = 4.207462
$ws.Cells.Item(14, 15).Value = 0.001822782814434402
$ws.Cells.Item(14, 16).Value = 0.001936017319896666
$ws.Cells.Item(14, 17).Value = 7.901238236890444
$ws.Cells.Item(14, 18).Value = 71.11114413201399
$ws.Cells.Item(14, 19).Value = 0.000359257599504591
$ws.Cells.Item(14, 20).Value = 0.0004107308728868184

$ws.Cells.Item(15, 7).Value = 5.633732333333334
$ws.Cells.Item(15, 8).Value = 16.901197
$ws.Cells.Item(15, 9).Value = 0.1970929266282699
$ws.Cells.Item(15, 10).Value = 0.2121524785267629
$ws.Cells.Item(15, 13).Value = 135.006546
$ws.Cells.Item(15, 14).Value = 270.013092
$ws.Cells.Item(15, 15).Value = 0.1754651225976237
$ws.Cells.Item(15, 16).Value = 0.1242435517446983
$ws.Cells.Item(15, 17).Value = 760.5907434118541
$ws.Cells.Item(15, 18).Value = 4563.544460471125
$ws.Cells.Item(15, 19).Value = 0.03458293453395384
$ws.Cells.Item(15, 20).Value = 0.02635857744360587

$ws.Cells.Item(16, 7).Value = 5.633732333333334
$ws.Cells.Item(16, 8).Value = 16.901197
$ws.Cells.Item(16, 9).Value = 0.1970929266282699
$ws.Cells.Item(16, 10).Value = 0.2121524785267629
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 0.6666666666666666
$ws.Cells.Item(16, 13).Value = 0.3895486666666667
$ws.Cells.Item(16, 14).Value = 1.168646
$ws.Cells.Item(16, 15).Value = 0.0005062880769826339
$ws.Cells.Item(16, 16).Value = 0.0005377395914277917
$ws.Cells.Item(16, 17).Value = 2.194612918806889
$ws.Cells.Item(16, 18).Value = 19.751516269262
$ws.Cells.Item(16, 19).Value = 0.00009978579880950615
$ws.Cells.Item(16, 20).Value = 0.0001140827871233748

$ws.Cells.Item(17, 7).Value = 6.087099
$ws.Cells.Item(17, 8).Value = 12.174198
$ws.Cells.Item(17, 9).Value = 0.2129537020222914
$ws.Cells.Item(17, 10).Value = 0.15281676675182
$ws.Cells.Item(17, 13).Value = 556.7425436666666
$ws.Cells.Item(17, 14).Value = 1670.227631
$ws.Cells.Item(17, 15).Value = 0.7235863858022448
$ws.Cells.Item(17, 16).Value = 0.7685368570853349
$ws.Cells.Item(17, 17).Value = 3388.946980810823
$ws.Cells.Item(17, 18).Value = 20333.68188486494
$ws.Cells.Item(17, 19).Value = 0.154090399589518
$ws.Cells.Item(17, 20).Value = 0.1174453176293864

$ws.Cells.Item(18, 7).Value = 6.087099
$ws.Cells.Item(18, 8).Value = 12.174198
$ws.Cells.Item(18, 9).Value = 0.2129537020222914
$ws.Cells.Item(18, 10).Value = 0.15281676675182
$ws.Cells.Item(18, 15).Value = 0.0986194207087145
$ws.Cells.Item(18, 16).Value = 0.1047458342586422
$ws.Cells.Item(18, 17).Value = 461.888165142249
$ws.Cells.Item(18, 18).Value = 2771.328990853494
$ws.Cells.Item(18, 19).Value = 0.02100137073121459
$ws.Cells.Item(18, 20).Value = 0.01600691972212771

$ws.Cells.Item(19, 7).Value = 6.087099
$ws.Cells.Item(19, 8).Value = 12.174198
$ws.Cells.Item(19, 9).Value = 0.2129537020222914
$ws.Cells.Item(19, 10).Value = 0.15281676675182
$ws.Cells.Item(19, 13).Value = 1.402487333333333
$ws.Cells.Item(19, 14).Value = 4.207462
$ws.Cells.Item(19, 15).Value = 0.001822782814434402
$ws.Cells.Item(19, 16).Value = 0.001936017319896666
$ws.Cells.Item(19, 17).Value = 8.537079244246
$ws.Cells.Item(19, 18).Value = 51.22247546547599
$ws.Cells.Item(19, 19).Value = 0.0003881683483164175
$ws.Cells.Item(19, 20).Value = 0.0002958559072021324

$ws.Cells.Item(20, 7).Value = 6.087099
$ws.Cells.Item(20, 8).Value = 12.174198
$ws.Cells.Item(20, 9).Value = 0.2129537020222914
$ws.Cells.Item(20, 10).Value = 0.15281676675182
$ws.Cells.Item(20, 13).Value = 135.006546
$ws.Cells.Item(20, 14).Value = 270.013092
$ws.Cells.Item(20, 15).Value = 0.1754651225976237
$ws.Cells.Item(20, 16).Value = 0.1242435517446983
$ws.Cells.Item(20, 17).Value = 821.7982111500542
$ws.Cells.Item(20, 18).Value = 3287.192844600217
$ws.Cells.Item(20, 19).Value = 0.03736594743295919
$ws.Cells.Item(20, 20).Value = 0.01898649786738724

$ws.Cells.Item(21, 7).Value = 6.087099
$ws.Cells.Item(21, 8).Value = 12.174198
$ws.Cells.Item(21, 9).Value = 0.2129537020222914
$ws.Cells.Item(21, 10).Value = 0.15281676675182
$ws.Cells.Item(21, 11).Value = 2
$ws.Cells.Item(21, 12).Value = 0.6666666666666666
$ws.Cells.Item(21, 13).Value = 0.3895486666666667
$ws.Cells.Item(21, 14).Value = 1.168646
$ws.Cells.Item(21, 15).Value = 0.0005062880769826339
$ws.Cells.Item(21, 16).Value = 0.0005377395914277917
$ws.Cells.Item(21, 17).Value = 2.371221299318
$ws.Cells.Item(21, 18).Value = 14.227327795908
$ws.Cells.Item(21, 19).Value = 0.0001078159202831988
$ws.Cells.Item(21, 20).Value = 0.00008217562571643982

$ws.Cells.Item(22, 7).Value = 7.094445666666666
$ws.Cells.Item(22, 8).Value = 21.283337
$ws.Cells.Item(22, 9).Value = 0.2481951531448182
$ws.Cells.Item(22, 10).Value = 0.2671593435583502
$ws.Cells.Item(22, 13).Value = 556.7425436666666
$ws.Cells.Item(22, 14).Value = 1670.227631
$ws.Cells.Item(22, 15).Value = 0.7235863858022448
$ws.Cells.Item(22, 16).Value = 0.7685368570853349
$ws.Cells.Item(22, 17).Value = 3949.77972636496
$ws.Cells.Item(22, 18).Value = 35548.01753728464
$ws.Cells.Item(22, 19).Value = 0.1795906338376937
$ws.Cells.Item(22, 20).Value = 0.2053218022393157

$ws.Cells.Item(23, 7).Value = 7.094445666666666
$ws.Cells.Item(23, 8).Value = 21.283337
$ws.Cells.Item(23, 9).Value = 0.2481951531448182
$ws.Cells.Item(23, 10).Value = 0.2671593435583502
$ws.Cells.Item(23, 15).Value = 0.0986194207087145
$ws.Cells.Item(23, 16).Value = 0.1047458342586422
$ws.Cells.Item(23, 17).Value = 538.3254801142623
$ws.Cells.Item(23, 18).Value = 4844.929321028361
$ws.Cells.Item(23, 19).Value = 0.02447686222585265
$ws.Cells.Item(23, 20).Value = 0.02798382832101059

$ws.Cells.Item(24, 7).Value = 7.094445666666666
$ws.Cells.Item(24, 8).Value = 21.283337
$ws.Cells.Item(24, 9).Value = 0.2481951531448182
$ws.Cells.Item(24, 10).Value = 0.2671593435583502
$ws.Cells.Item(24, 13).Value = 1.402487333333333
$ws.Cells.Item(24, 14).Value = 4.207462
$ws.Cells.Item(24, 15).Value = 0.001822782814434402
$ws.Cells.Item(24, 16).Value = 0.001936017319896666
$ws.Cells.Item(24, 17).Value = 9.949870184521554
$ws.Cells.Item(24, 18).Value = 89.54883166069399
$ws.Cells.Item(24, 19).Value = 0.0004524058597782893
$ws.Cells.Item(24, 20).Value = 0.0005172251163011897

$ws.Cells.Item(25, 7).Value = 7.094445666666666
$ws.Cells.Item(25, 8).Value = 21.283337
$ws.Cells.Item(25, 9).Value = 0.2481951531448182
$ws.Cells.Item(25, 10).Value = 0.2671593435583502
$ws.Cells.Item(25, 13).Value = 135.006546
$ws.Cells.Item(25, 14).Value = 270.013092
$ws.Cells.Item(25, 15).Value = 0.1754651225976237
$ws.Cells.Item(25, 16).Value = 0.1242435517446983
$ws.Cells.Item(25, 17).Value = 957.7966052413341
$ws.Cells.Item(25, 18).Value = 5746.779631448005
$ws.Cells.Item(25, 19).Value = 0.04354959297469151
$ws.Cells.Item(25, 20).Value = 0.03319282572547153

$ws.Cells.Item(26, 7).Value = 7.094445666666666
$ws.Cells.Item(26, 8).Value = 21.283337
$ws.Cells.Item(26, 9).Value = 0.2481951531448182
$ws.Cells.Item(26, 10).Value = 0.2671593435583502
$ws.Cells.Item(26, 11).Value = 2
$ws.Cells.Item(26, 12).Value = 0.6666666666666666
$ws.Cells.Item(26, 13).Value = 0.3895486666666667
$ws.Cells.Item(26, 14).Value = 1.168646
$ws.Cells.Item(26, 15).Value = 0.0005062880769826339
$ws.Cells.Item(26, 16).Value = 0.0005377395914277917
$ws.Cells.Item(26, 17).Value = 2.763631850189111
$ws.Cells.Item(26, 18).Value = 2.763631850189111
$ws.Cells.Item(26, 19).Value = 0.0001256582468021004
$ws.Cells.Item(26, 20).Value = 0.0001436621562511843

